# Auto-generated edit script applying the Zodiark_Profits price-refresh diff
# across the 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 95.916664
$ws.Range("J12").Value = 75.28570999999999
$ws.Range("L12").Value = 75.28570999999999
$ws.Range("N12").Value = -415.28571
$ws.Range("H21").Value = 9935.286
$ws.Range("I21").Value = 10806.2
$ws.Range("J21").Value = 7758
$ws.Range("K21").Value = 10806.2
$ws.Range("L21").Value = 7758
$ws.Range("M21").Value = -10338.2
$ws.Range("N21").Value = -8694
$ws.Range("H23").Value = 9935.286
$ws.Range("I23").Value = 10806.2
$ws.Range("J23").Value = 7758
$ws.Range("K23").Value = 10806.2
$ws.Range("L23").Value = 7758
$ws.Range("M23").Value = -10572.2
$ws.Range("N23").Value = -8226
$ws.Range("H28").Value = 774.2143
$ws.Range("I28").Value = 1031.3334
$ws.Range("K28").Value = 1031.3334
$ws.Range("M28").Value = -546.3334
$ws.Range("H62").Value = 4327.609
$ws.Range("I62").Value = 4340.3335
$ws.Range("J62").Value = 4194
$ws.Range("K62").Value = 4340.3335
$ws.Range("L62").Value = 4194
$ws.Range("M62").Value = -3716.3335
$ws.Range("N62").Value = -5442
$ws.Range("H65").Value = 4327.609
$ws.Range("I65").Value = 4340.3335
$ws.Range("J65").Value = 4194
$ws.Range("K65").Value = 21701.6675
$ws.Range("L65").Value = 20970
$ws.Range("M65").Value = -18581.6675
$ws.Range("N65").Value = -27210
$ws.Range("H98").Value = 1994.1875
$ws.Range("I98").Value = 1993.5714
$ws.Range("K98").Value = 1993.5714
$ws.Range("M98").Value = -495.5714
$ws.Range("H122").Value = 1994.1875
$ws.Range("I122").Value = 1993.5714
$ws.Range("K122").Value = 5980.7142
$ws.Range("M122").Value = -3530.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5608.1772
$ws.Range("I32").Value = 3972.9714
$ws.Range("K32").Value = 3972.9714
$ws.Range("M32").Value = -3685.9714
$ws.Range("H45").Value = 1518.7931
$ws.Range("I45").Value = 1383.95
$ws.Range("J45").Value = 1818.4445
$ws.Range("K45").Value = 1383.95
$ws.Range("L45").Value = 1818.4445
$ws.Range("M45").Value = -1006.95
$ws.Range("N45").Value = -2572.4445
$ws.Range("H61").Value = 2606.1538
$ws.Range("I61").Value = 2378.818
$ws.Range("K61").Value = 2378.818
$ws.Range("M61").Value = -2166.818
$ws.Range("H102").Value = 125387500
$ws.Range("I102").Value = 167016670
$ws.Range("J102").Value = 500000
$ws.Range("K102").Value = 167016670
$ws.Range("L102").Value = 500000
$ws.Range("M102").Value = -167015048
$ws.Range("N102").Value = -503244
$ws.Range("H132").Value = 6739.047
$ws.Range("I132").Value = 4546.875
$ws.Range("J132").Value = 22084.25
$ws.Range("K132").Value = 13640.625
$ws.Range("L132").Value = 66252.75
$ws.Range("M132").Value = -11110.625
$ws.Range("N132").Value = -71312.75
$ws.Range("H136").Value = 2606.1538
$ws.Range("I136").Value = 2378.818
$ws.Range("K136").Value = 7136.454000000001
$ws.Range("M136").Value = -4586.454000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1409.7843
$ws.Range("I134").Value = 1464.8605
$ws.Range("K134").Value = 4394.5815
$ws.Range("M134").Value = -1859.5815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 229433.16
$ws.Range("I31").Value = 3768.7256
$ws.Range("K31").Value = 3768.7256
$ws.Range("M31").Value = -3473.7256
$ws.Range("H34").Value = 229433.16
$ws.Range("I34").Value = 3768.7256
$ws.Range("K34").Value = 3768.7256
$ws.Range("M34").Value = -3566.7256

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 9999
$ws.Range("J48").Value = 9999
$ws.Range("L48").Value = 29997
$ws.Range("N48").Value = -30497
$ws.Range("H68").Value = 1399
$ws.Range("J68").Value = 1900
$ws.Range("L68").Value = 5700
$ws.Range("N68").Value = -7322
$ws.Range("H71").Value = 1399
$ws.Range("J71").Value = 1900
$ws.Range("L71").Value = 17100
$ws.Range("N71").Value = -25212
$ws.Range("H81").Value = 6210.5
$ws.Range("I81").Value = 1057.8
$ws.Range("K81").Value = 3173.4
$ws.Range("M81").Value = -2050.4
$ws.Range("H84").Value = 6210.5
$ws.Range("I84").Value = 1057.8
$ws.Range("K84").Value = 9520.199999999999
$ws.Range("M84").Value = -3904.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4101.4
$ws.Range("I80").Value = 4003
$ws.Range("J80").Value = 4249
$ws.Range("K80").Value = 4003
$ws.Range("L80").Value = 4249
$ws.Range("M80").Value = -3005
$ws.Range("N80").Value = -6245
$ws.Range("H83").Value = 4101.4
$ws.Range("I83").Value = 4003
$ws.Range("J83").Value = 4249
$ws.Range("K83").Value = 20015
$ws.Range("L83").Value = 21245
$ws.Range("M83").Value = -15023
$ws.Range("N83").Value = -31229
$ws.Range("H102").Value = 1197.1298
$ws.Range("I102").Value = 1181.1262
$ws.Range("K102").Value = 1181.1262
$ws.Range("M102").Value = 440.8738000000001
$ws.Range("H122").Value = 3357
$ws.Range("I122").Value = 1628.3334
$ws.Range("K122").Value = 4885.0002
$ws.Range("M122").Value = -2435.0002
$ws.Range("H132").Value = 7400.4165
$ws.Range("I132").Value = 6814.722
$ws.Range("K132").Value = 20444.166
$ws.Range("M132").Value = -17914.166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 43486300
$ws.Range("I7").Value = 125004450
$ws.Range("J7").Value = 9955.267
$ws.Range("K7").Value = 125004450
$ws.Range("L7").Value = 9955.267
$ws.Range("M7").Value = -125004338
$ws.Range("N7").Value = -10179.267
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H126").Value = 43486300
$ws.Range("I126").Value = 125004450
$ws.Range("J126").Value = 9955.267
$ws.Range("K126").Value = 375013350
$ws.Range("L126").Value = 29865.801
$ws.Range("M126").Value = -375010880
$ws.Range("N126").Value = -34805.801
$ws.Range("H132").Value = 5016.6523
$ws.Range("I132").Value = 4569.15
$ws.Range("K132").Value = 13707.45
$ws.Range("M132").Value = -11177.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4632020.5
$ws.Range("I96").Value = 9260681
$ws.Range("J96").Value = 3360
$ws.Range("K96").Value = 9260681
$ws.Range("L96").Value = 3360
$ws.Range("M96").Value = -9259308
$ws.Range("N96").Value = -6106
$ws.Range("H132").Value = 2051.0469
$ws.Range("J132").Value = 2515.0312
$ws.Range("L132").Value = 7545.0936
$ws.Range("N132").Value = -12605.0936
